# Adds the Dec 21 ("day" 21) calendar entry to row 22 of Tabelle1 and
# moves the active selection to E23, matching the upstream commit
# "Add files via upload".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Write order matches the shared-string append order used by the
# original commit: title, person, image_url, then text.
$ws.Range("B22").Value = "Energiestaubsauger"
$ws.Range("D22").Value = "Energie zentrieren"
$ws.Range("E22").Value = "https://us.123rf.com/450wm/dxinerz/dxinerz1801/dxinerz180100321/92601781-yoga-aerobic-bewegung-fitness-sport-symbolbild-kann-auch-f%C3%BCr-fitness-erholung-verwendet-werden.jpg?ver=6"
$ws.Range("C22").Value = "Am Sonntag war Padeltraining. Ich bin endlich wieder mal gegangen und spielte zuerst gleich mit jemanden, den ich gut kenne. Aktuell ist er leider sehr negativ eingestellt gegenüber allem und kommt gefühlt jedes Mal direkt von der Party und ist dementsprechend fit und munter beim Spielen. Ich habe mich über diese negative Energie aufgeregt. Früher habe ich mich oft reinziehen lassen. Dieses Mal habe ich mir danach bewusst eine Pause gegönnt und den Energiestaubsauger geübt, um mich bewusst abzugrenzen von diesen Energien und dass sie nicht zu meinen werden. Es hat geholfen, aber ich darf noch weiter üben. Und ich habe gemerkt, dass ich aktuell gerade Mühe habe mit dieser negativen Energe ""alles ist scheisse"". Davon darf ich mich noch etwas mehr abgrenzen."

[void]$ws.Range("E23").Select()
